$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4, 2).Value = 1698722
$ws.Cells.Item(4, 3).Value = 12286
$ws.Cells.Item(4, 4).Value = 457232
$ws.Cells.Item(4, 5).Value = 1141842
$ws.Cells.Item(4, 7).Value = 348
$ws.Cells.Item(4, 8).Value = 99648

# --- Row 7: Espana ---
$ws.Cells.Item(7, 2).Value = 282480
$ws.Cells.Item(7, 5).Value = 58685
$ws.Cells.Item(7, 8).Value = 26837

# --- Row 16: Canada ---
$ws.Cells.Item(16, 2).Value = 85677
$ws.Cells.Item(16, 3).Value = 978
$ws.Cells.Item(16, 5).Value = 34932
$ws.Cells.Item(16, 7).Value = 114
$ws.Cells.Item(16, 8).Value = 6538

# --- Rows 34/35: Sudafrica overtakes Indonesia in ranking ---
# Row 34 becomes Sudafrica with new (updated) data
$ws.Cells.Item(34, 1).Value = "Sudafrica"
$ws.Cells.Item(34, 2).Value = 23615
$ws.Cells.Item(34, 3).Value = 1032
$ws.Cells.Item(34, 4).Value = 11917
$ws.Cells.Item(34, 5).Value = 11217
$ws.Cells.Item(34, 7).Value = 52
$ws.Cells.Item(34, 8).Value = 481

# Row 35 becomes Indonesia, keeping its previous (unchanged) data
$ws.Cells.Item(35, 1).Value = "Indonesia"
$ws.Cells.Item(35, 2).Value = 22750
$ws.Cells.Item(35, 3).Value = 479
$ws.Cells.Item(35, 4).Value = 5642
$ws.Cells.Item(35, 5).Value = 15717
$ws.Cells.Item(35, 7).Value = 19
$ws.Cells.Item(35, 8).Value = 1391

# --- Rows 142-145: Ruanda moves up in ranking, others shift down one row ---
# Row 142 becomes Ruanda with new (updated) data
$ws.Cells.Item(142, 1).Value = "Ruanda"
$ws.Cells.Item(142, 2).Value = 336
$ws.Cells.Item(142, 3).Value = 9
$ws.Cells.Item(142, 4).Value = 238
$ws.Cells.Item(142, 5).Value = 98
$ws.Cells.Item(142, 7).Value = 0
$ws.Cells.Item(142, 8).Value = 0

# Row 143 becomes Isla de Man (former row 142 data)
$ws.Cells.Item(143, 1).Value = "Isla de Man"
$ws.Cells.Item(143, 2).Value = 336
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 303
$ws.Cells.Item(143, 5).Value = 9
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 24

# Row 144 becomes Mauricio (former row 143 data)
$ws.Cells.Item(144, 1).Value = "Mauricio"
$ws.Cells.Item(144, 2).Value = 334
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 322
$ws.Cells.Item(144, 5).Value = 2
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 10

# Row 145 becomes Guayana Francesa (former row 144 data)
$ws.Cells.Item(145, 1).Value = "Guayana Francesa"
$ws.Cells.Item(145, 2).Value = 328
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 145
$ws.Cells.Item(145, 5).Value = 182
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 1

# --- Update timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 21:35"
